$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.023267232260603
$ws.Cells.Item(2, 4).Value = 1.029146288686311
$ws.Cells.Item(2, 5).Value = 1.02391312102178
$ws.Cells.Item(2, 9).Value = 1.0325998660753
$ws.Cells.Item(2, 10).Value = 1.028448410760337
$ws.Cells.Item(2, 11).Value = 1.031961047169429
$ws.Cells.Item(2, 12).Value = 1.026743163461279
$ws.Cells.Item(2, 14).Value = 1.013563638581738

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.023994022218411
$ws.Cells.Item(3, 4).Value = 1.029686919012857
$ws.Cells.Item(3, 5).Value = 1.024523839692532
$ws.Cells.Item(3, 9).Value = 1.032744194110722
$ws.Cells.Item(3, 10).Value = 1.028814760757309
$ws.Cells.Item(3, 11).Value = 1.032310482425216
$ws.Cells.Item(3, 12).Value = 1.027161396635822
$ws.Cells.Item(3, 14).Value = 1.01368408856369

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.024464740642486
$ws.Cells.Item(4, 4).Value = 1.030036958671803
$ws.Cells.Item(4, 5).Value = 1.024919797738655
$ws.Cells.Item(4, 9).Value = 1.032836348136384
$ws.Cells.Item(4, 10).Value = 1.02905155962242
$ws.Cells.Item(4, 11).Value = 1.032536097416402
$ws.Cells.Item(4, 12).Value = 1.027432109903117
$ws.Cells.Item(4, 14).Value = 1.013761941864797

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.0246627331913
$ws.Cells.Item(5, 4).Value = 1.030184164952103
$ws.Cells.Item(5, 5).Value = 1.025086443715861
$ws.Cells.Item(5, 9).Value = 1.032874792885387
$ws.Cells.Item(5, 10).Value = 1.029151047654226
$ws.Cells.Item(5, 11).Value = 1.032630826472655
$ws.Cells.Item(5, 12).Value = 1.0275459374186
$ws.Cells.Item(5, 14).Value = 1.013794650364313

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.024695982964449
$ws.Cells.Item(6, 4).Value = 1.030208884384227
$ws.Cells.Item(6, 5).Value = 1.025114435109993
$ws.Cells.Item(6, 9).Value = 1.032881230496011
$ws.Cells.Item(6, 10).Value = 1.029167748444152
$ws.Cells.Item(6, 11).Value = 1.032646724848662
$ws.Cells.Item(6, 12).Value = 1.027565050667965
$ws.Cells.Item(6, 14).Value = 1.013800141018636

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.024467385828317
$ws.Cells.Item(7, 4).Value = 1.030038925456962
$ws.Cells.Item(7, 5).Value = 1.024922023745183
$ws.Cells.Item(7, 9).Value = 1.032836863004408
$ws.Cells.Item(7, 10).Value = 1.029052889232559
$ws.Cells.Item(7, 11).Value = 1.032537363662031
$ws.Cells.Item(7, 12).Value = 1.027433630796617
$ws.Cells.Item(7, 14).Value = 1.013762379000582

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.023512762684099
$ws.Cells.Item(8, 4).Value = 1.02932895135941
$ws.Cells.Item(8, 5).Value = 1.024119353057959
$ws.Cells.Item(8, 9).Value = 1.032648897742623
$ws.Cells.Item(8, 10).Value = 1.028572272234907
$ws.Cells.Item(8, 11).Value = 1.032079241646674
$ws.Cells.Item(8, 12).Value = 1.026884488100127
$ws.Cells.Item(8, 14).Value = 1.013604362721834

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.021834029181753
$ws.Cells.Item(9, 4).Value = 1.028079627688584
$ws.Cells.Item(9, 5).Value = 1.022711026300132
$ws.Cells.Item(9, 9).Value = 1.032308254991301
$ws.Cells.Item(9, 10).Value = 1.027723478821064
$ws.Cells.Item(9, 11).Value = 1.031268261077129
$ws.Cells.Item(9, 12).Value = 1.025917572666616
$ws.Cells.Item(9, 14).Value = 1.013325280577399

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.020717303703897
$ws.Cells.Item(10, 4).Value = 1.027248037023614
$ws.Cells.Item(10, 5).Value = 1.021776353228408
$ws.Cells.Item(10, 9).Value = 1.032074879435337
$ws.Cells.Item(10, 10).Value = 1.027156431780957
$ws.Cells.Item(10, 11).Value = 1.030725199557231
$ws.Cells.Item(10, 12).Value = 1.025273553033571
$ws.Cells.Item(10, 14).Value = 1.013138825106761

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.020234351244219
$ws.Cells.Item(11, 4).Value = 1.026888281103367
$ws.Cells.Item(11, 5).Value = 1.021372654354499
$ws.Cells.Item(11, 9).Value = 1.031972347380985
$ws.Cells.Item(11, 10).Value = 1.026910630592255
$ws.Cells.Item(11, 11).Value = 1.030489495284539
$ws.Cells.Item(11, 12).Value = 1.024994844645536
$ws.Cells.Item(11, 14).Value = 1.013057998657793

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.020055052967709
$ws.Cells.Item(12, 4).Value = 1.026754703147339
$ws.Cells.Item(12, 5).Value = 1.021222858093924
$ws.Cells.Item(12, 9).Value = 1.031934041157772
$ws.Cells.Item(12, 10).Value = 1.026819290521181
$ws.Cells.Item(12, 11).Value = 1.030401862322283
$ws.Cells.Item(12, 12).Value = 1.024891344985301
$ws.Cells.Item(12, 14).Value = 1.013027963068904

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.020093508870612
$ws.Cells.Item(13, 4).Value = 1.026783353712513
$ws.Cells.Item(13, 5).Value = 1.021254982830697
$ws.Cells.Item(13, 9).Value = 1.03194226796985
$ws.Cells.Item(13, 10).Value = 1.026838884999045
$ws.Cells.Item(13, 11).Value = 1.030420663572466
$ws.Cells.Item(13, 12).Value = 1.024913544855552
$ws.Cells.Item(13, 14).Value = 1.013034406387475

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.020219528508685
$ws.Cells.Item(14, 4).Value = 1.026877238447059
$ws.Cells.Item(14, 5).Value = 1.021360268964961
$ws.Cells.Item(14, 9).Value = 1.031969185484715
$ws.Cells.Item(14, 10).Value = 1.026903081172064
$ws.Cells.Item(14, 11).Value = 1.030482253178689
$ws.Cells.Item(14, 12).Value = 1.024986288811468
$ws.Cells.Item(14, 14).Value = 1.01305551617006

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.02029718558929
$ws.Cells.Item(15, 4).Value = 1.026935090805247
$ws.Cells.Item(15, 5).Value = 1.021425159892586
$ws.Cells.Item(15, 9).Value = 1.03198574096379
$ws.Cells.Item(15, 10).Value = 1.026942629484844
$ws.Cells.Item(15, 11).Value = 1.030520189759652
$ws.Cells.Item(15, 12).Value = 1.025031112128581
$ws.Cells.Item(15, 14).Value = 1.013068520889146

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.020749368423943
$ws.Cells.Item(16, 4).Value = 1.027271919934673
$ws.Cells.Item(16, 5).Value = 1.021803167080106
$ws.Cells.Item(16, 9).Value = 1.032081653071226
$ws.Cells.Item(16, 10).Value = 1.027172739310788
$ws.Cells.Item(16, 11).Value = 1.03074083090648
$ws.Cells.Item(16, 12).Value = 1.025292053418412
$ws.Cells.Item(16, 14).Value = 1.013144187436324

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.021033172077654
$ws.Cells.Item(17, 4).Value = 1.027483293323798
$ws.Cells.Item(17, 5).Value = 1.02204055578164
$ws.Cells.Item(17, 9).Value = 1.032141420919337
$ws.Cells.Item(17, 10).Value = 1.027317010884825
$ws.Cells.Item(17, 11).Value = 1.03087908571841
$ws.Cells.Item(17, 12).Value = 1.025455777955861
$ws.Cells.Item(17, 14).Value = 1.01319162730215

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.021198767414307
$ws.Cells.Item(18, 4).Value = 1.027606615384795
$ws.Cells.Item(18, 5).Value = 1.022179118939021
$ws.Cells.Item(18, 9).Value = 1.03217613973905
$ws.Cells.Item(18, 10).Value = 1.027401136286001
$ws.Cells.Item(18, 11).Value = 1.030959673757083
$ws.Cells.Item(18, 12).Value = 1.025551290691413
$ws.Cells.Item(18, 14).Value = 1.013219289457844

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.021255240829085
$ws.Cells.Item(19, 4).Value = 1.02764867030678
$ws.Cells.Item(19, 5).Value = 1.022226382011099
$ws.Cells.Item(19, 9).Value = 1.032187953717195
$ws.Cells.Item(19, 10).Value = 1.027429816460317
$ws.Cells.Item(19, 11).Value = 1.030987143049661
$ws.Cells.Item(19, 12).Value = 1.025583860560014
$ws.Cells.Item(19, 14).Value = 1.01322872004534

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.021002716672475
$ws.Cells.Item(20, 4).Value = 1.027460611682273
$ws.Cells.Item(20, 5).Value = 1.022015076033092
$ws.Cells.Item(20, 9).Value = 1.032135023153938
$ws.Cells.Item(20, 10).Value = 1.027301534565341
$ws.Cells.Item(20, 11).Value = 1.030864257824462
$ws.Cells.Item(20, 12).Value = 1.025438210301516
$ws.Cells.Item(20, 14).Value = 1.013186538351567

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.020182416316287
$ws.Cells.Item(21, 4).Value = 1.02684959029686
$ws.Cells.Item(21, 5).Value = 1.021329260513703
$ws.Cells.Item(21, 9).Value = 1.031961265047652
$ws.Cells.Item(21, 10).Value = 1.026884178052441
$ws.Cells.Item(21, 11).Value = 1.030464118823654
$ws.Cells.Item(21, 12).Value = 1.024964866825744
$ws.Cells.Item(21, 14).Value = 1.013049300221591

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.019667192904561
$ws.Cells.Item(22, 4).Value = 1.026465715374399
$ws.Cells.Item(22, 5).Value = 1.020898961657924
$ws.Cells.Item(22, 9).Value = 1.031850736790243
$ws.Cells.Item(22, 10).Value = 1.026621547016602
$ws.Cells.Item(22, 11).Value = 1.030212062434998
$ws.Cells.Item(22, 12).Value = 1.024667402661271
$ws.Cells.Item(22, 14).Value = 1.012962937910681

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.019940271466621
$ws.Cells.Item(23, 4).Value = 1.026669185748111
$ws.Cells.Item(23, 5).Value = 1.021126985151321
$ws.Cells.Item(23, 9).Value = 1.031909450893921
$ws.Cells.Item(23, 10).Value = 1.026760793311821
$ws.Cells.Item(23, 11).Value = 1.030345726645402
$ws.Cells.Item(23, 12).Value = 1.024825079718262
$ws.Cells.Item(23, 14).Value = 1.013008727180594

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.021016477981847
$ws.Cells.Item(24, 4).Value = 1.027470860442587
$ws.Cells.Item(24, 5).Value = 1.02202658893112
$ws.Cells.Item(24, 9).Value = 1.032137914469942
$ws.Cells.Item(24, 10).Value = 1.027308527727748
$ws.Cells.Item(24, 11).Value = 1.030870958077457
$ws.Cells.Item(24, 12).Value = 1.025446148322237
$ws.Cells.Item(24, 14).Value = 1.013188837856276

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.022267602550434
$ws.Cells.Item(25, 4).Value = 1.028402388950305
$ws.Cells.Item(25, 5).Value = 1.023074379180041
$ws.Cells.Item(25, 9).Value = 1.032397430346409
$ws.Cells.Item(25, 10).Value = 1.027943127115876
$ws.Cells.Item(25, 11).Value = 1.031478349763246
$ws.Cells.Item(25, 12).Value = 1.026167445974109
$ws.Cells.Item(25, 14).Value = 1.013397502656659
